$wb = $excel.ActiveWorkbook

# zh-cn sheet: update Correspond Handoff/Handback datetimes
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-23 07:27:27"
$wsZh.Range("E3").Value = "2016-03-23 07:27:27"
$wsZh.Range("H2").Value = "2016-03-23 07:28:08"
$wsZh.Range("H3").Value = "2016-03-23 07:28:08"

# de-de sheet: update Correspond Handoff/Handback datetimes
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-23 07:27:35"
$wsDe.Range("E3").Value = "2016-03-23 07:27:35"
$wsDe.Range("H2").Value = "2016-03-23 07:28:21"
$wsDe.Range("H3").Value = "2016-03-23 07:28:21"
